# Fixed alignment in one of the header boxes
#
# This script reproduces four logical changes on slide 1:
#   1. Drop the stray <a:endParaRPr> left behind on the first paragraph of
#      the "magic" textbox (paragraph mark formatting is no longer needed
#      once the paragraph content run is rewritten).
#   2. Same cleanup for the "version" textbox.
#   3. Merge the two runs that spell out "36-byte header (first 20 bytes)"
#      into a single run (the leading "36" and the rest were split into two
#      <a:r> runs with identical formatting - collapse them into one).
#   4. Same run-merge for "36-byte header (next 16 bytes)".
#   5/6. Nudge the "gltfFormat" box and the trailing "..." box up very
#      slightly (a few EMU) to fix their vertical alignment with the row.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- helper: EMU -> points, nudged so the engine's point->EMU rounding
# lands back on the exact integer EMU we want (avoids off-by-one-EMU
# truncation when going through the float Top/Left properties). ---
function EmuToPt([double]$emu) {
    return ($emu / 12700.0) + 0.00001
}

# 1) "magic" textbox: remove the leftover endParaRPr on paragraph 1 without
#    touching paragraph 2 (which keeps its own runs/colors/endParaRPr as-is).
#    Deleting just paragraph 1 (text + its paragraph mark) and retyping the
#    text back in front of paragraph 2 rebuilds paragraph 1 with a fresh run
#    and no endParaRPr, while paragraph 2 is left completely untouched.
$shpMagic = $s.Shapes.Item("TextBox 3")
$trMagic = $shpMagic.TextFrame.TextRange
$trMagic.Paragraphs(1).Delete()
[void]$trMagic.InsertBefore("magic" + [char]13)

# 2) "version" textbox: identical fix.
$shpVersion = $s.Shapes.Item("TextBox 4")
$trVersion = $shpVersion.TextFrame.TextRange
$trVersion.Paragraphs(1).Delete()
[void]$trVersion.InsertBefore("version" + [char]13)

# 3) "36-byte header (first 20 bytes)" textbox: merge the two runs into one.
#    Re-assigning identical text is a no-op for this engine, so first set a
#    throwaway value, then set the final text - this forces the paragraph's
#    runs to be rebuilt as a single run while preserving the trailing
#    endParaRPr untouched.
$shpHeader1 = $s.Shapes.Item("TextBox 13")
$trHeader1 = $shpHeader1.TextFrame.TextRange
$trHeader1.Text = "x"
$trHeader1.Text = "36-byte header (first 20 bytes)"

# 4) "36-byte header (next 16 bytes)" textbox: same run-merge fix.
$shpHeader2 = $s.Shapes.Item("TextBox 36")
$trHeader2 = $shpHeader2.TextFrame.TextRange
$trHeader2.Text = "x"
$trHeader2.Text = "36-byte header (next 16 bytes)"

# 5) "gltfFormat" box: move up from y=1862283 EMU to y=1858875 EMU.
$shpGltf = $s.Shapes.Item("TextBox 15")
$shpGltf.Top = EmuToPt 1858875

# 6) trailing "..." box: move up from y=1911349 EMU to y=1907941 EMU.
$shpEllipsis = $s.Shapes.Item("TextBox 40")
$shpEllipsis.Top = EmuToPt 1907941
